# render website, remove theme (not needed) from docs
#
# 1. Add a new "Abstract Title" paragraph style (based on Normal, followed
#    by Abstract) with centered, bold, small, blue text.
# 2. Abstract style: reduce space-before from 300 (15pt) to 100 (5pt) twips.
# 3. Add a new "Footnote Block Text" paragraph style (based on Footnote
#    Text, followed by Footnote Text) with block-quote-like indentation.

$d = $word.ActiveDocument

# --- 1. "Abstract Title" style -------------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = $d.Styles.Item("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles.Item("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# --- 2. "Abstract" style: before-spacing 300 -> 100 -----------------------
$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. "Footnote Block Text" style ---------------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = $d.Styles.Item("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles.Item("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "styles updated"
